$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the refreshed "time_taken" timestamps on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:19:13.369939"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:13.369947"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:13.369950"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:13.369953"
$dataSheet.Range("F6").Value = "2021-10-05 14:19:13.369956"

# --- Add the new "metadata" sheet, placed right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used by "data" (left/right 0.75in, top/bottom 1in,
# header/footer 0.5in -> 54/54/72/72/36/36 points).
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Bring the header style (bold font, thin border, centered/top aligned -
# the style already used by data!B1:F1) over onto the new header row.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Auditory Neuropathy Spectrum Disorde"
$metaSheet.Range("C2").Value = 260
$metaSheet.Range("E2").Value = "2019-02-17T23:29:21.469051Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:13.366554"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/260/?format=json"

# "1.8" must land as literal text, not be coerced to the number 1.8, and
# must keep the sheet's default (unstyled) cell formatting - stage it in a
# throwaway cell formatted as Text, copy only the *value* over, then wipe
# the scratch cell (and its formatting) away again.
$metaSheet.Range("Z1").NumberFormat = "@"
$metaSheet.Range("Z1").Value = "1.8"
$metaSheet.Range("Z1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)     # xlPasteValues
$metaSheet.Range("Z1").Clear()
$excel.CutCopyMode = $false

# Restore "data" as the active sheet/tab (unchanged by this edit).
$dataSheet.Activate()

Write-Output "metadata sheet added; timestamps refreshed"
